$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A and C (Excel's ColumnWidth setter adds a 5/6 padding offset
# internally in this runtime, so we subtract it to land on the exact target
# widths stored in the saved XML: 94 and 785).
$ws.Columns.Item(1).ColumnWidth = 94 - 5/6
$ws.Columns.Item(3).ColumnWidth = 785 - 5/6

# Append the three new Q&A rows (198-200) to the sheet.
$ws.Cells.Item(198, 1).Value = "Why can't I add 251 curve shades to my log?"
$ws.Cells.Item(198, 2).Value = "llama3.2:latest"
$ws.Cells.Item(198, 3).Value = "You cannot add 251 curve shades to your log because the maximum number of curve shades per plot is 250."

$ws.Cells.Item(199, 1).Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Cells.Item(199, 2).Value = "llama3.2:latest"
$row199c = @"
To apply a curve shading with the specified name, follow these steps:
1. On the Curve tab, select Curve Shading.
2. The Curve Shading Information dialog box will open.
3. From Curve Shading, select an unused curve number.
4. Select Display Curve Shadeto display the shading on the plot.
5. Enter a curve shadingName as "Hydrocarbon bearing zone highlighted".
6. Select the Track where the GR curve displays.
7. In To select <fixed curve value>, enter 40.
8. In the Zone Attributes grid, Type column, select Pattern.
9. In the Pattern column, select Sandstone.
10. Note: Sandstone.vob must be present in your ODF file as a lithology.
This will apply the specified curve shading with the name "Hydrocarbon bearing zone highlighted" to the GR curve when its value is less than 40 API.
"@
$ws.Cells.Item(199, 3).Value = $row199c
# The embedded line breaks make the runtime auto-expand the row height on
# save; AutoFit() (with WrapText left off) restores the standard 15pt row
# height and avoids emitting a stray ht/customHeight attribute on the row.
$ws.Rows.Item(199).AutoFit()

$ws.Cells.Item(200, 1).Value = "I have 20000 modifiers added ty log, why I can't I add anymore?"
$ws.Cells.Item(200, 2).Value = "llama3.2:latest"
$ws.Cells.Item(200, 3).Value = 'You cannot add more than 20000 modifiers per plot because the system has a limit on the number of modifiers that can be applied to a single well log. This is specified in the GEO application''s documentation under "Modifiers" section, which states that there are 20000 allowed modifiers per plot.'
